# Update crypto price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.034.41"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.827.38"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.75"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4339"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3679"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07282"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8459"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.70"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.832.87"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.668"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07069"
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.298"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.54"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008783"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.93"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.110.01"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.149"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.053.27"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.989"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.56"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.219"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.31"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.236"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.00"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08754"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("E32").Value = "  -2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7416"
$ws.Range("E33").Value = "  -3.82%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.908"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.440"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05249"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.233"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.875"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1703"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.588"
$ws.Range("E44").Value = "  -1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.62"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4772"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.13"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.936"
$ws.Range("E48").Value = "  +5.35%  "
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.663"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06337"
